# Fix choice name REGEX and associated test
#
# The "choices" sheet (3rd sheet) gets a new row demonstrating a value that
# is "skipped" by the (fixed) choice-name regex: a cell that only contains a
# leading single-quote (stored by Excel as an empty, quote-prefixed cell)
# followed by a cell with the literal text "skipped".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Row 7: new shared string "skipped"
$ws.Range("A7").Value = "skipped"

# Row 6: a cell whose only content is a quote-prefix (leading apostrophe),
# which Excel stores as an empty cell carrying a quotePrefix style.
$ws.Range("A6").Value = "'"
$ws.Range("A6").Value = ""

# Move the active selection down to A8, past the newly added rows.
$null = $ws.Range("A8").Select()
